# Add new columns I ("I0") and J ("IF") to the stats table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: set values, then copy the existing header style (H1) onto
# the two new header cells so they pick up the same bold/border/centered
# formatting (cellXfs index 1) as the rest of row 1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 4

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 3

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 4
